# MSFT, MU, OXY, XOM
# Update the long-term (terminal) growth-rate assumption on the Model sheet
# from 3% to 2%. Every downstream projection, NPV, and per-share figure in
# columns AU:DM (row 16) plus AQ22/AQ23/AQ25 is formula-driven off this one
# input cell, so changing it alone ripples through the whole model.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model")
$ws.Activate()

$ws.Range("AQ20").Value = 0.02

# Reflect the author's final cell selection on the Model sheet (cosmetic
# view-state, matches the saved workbook's active cell/selection).
$ws.Range("AR24").Select()
